$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update title/header shared strings ---
$ws.Range("A8").Value = "Volume 32   Number  8"
$ws.Range("C9").Value = "Report Covering the Week  2/17/2025  Through  2/23/2025"

# --- Update crime statistics table (rows 14-31) ---
# Row 14
$ws.Range("C14").Copy($ws.Range("F14"))
# Row 15
$ws.Range("D15").Value = 1
$ws.Range("E15").Value = -100
$ws.Range("F15").Value = 2
$ws.Range("G15").Value = 2
$ws.Range("H15").Value = 0
$ws.Range("J15").Value = 4
$ws.Range("K15").Value = 50
$ws.Range("N15").Value = -14.285714285714
# Row 16
$ws.Range("C16").Value = 9
$ws.Range("D16").Value = 7
$ws.Range("E16").Value = 28.571428571428
$ws.Range("F16").Value = 20
$ws.Range("G16").Value = 26
$ws.Range("H16").Value = -23.076923076923
$ws.Range("I16").Value = 33
$ws.Range("J16").Value = 52
$ws.Range("K16").Value = -36.538461538461
$ws.Range("L16").Value = -8.333333333333
$ws.Range("M16").Value = -34
$ws.Range("N16").Value = -83.743842364532
# Row 17
$ws.Range("D17").Value = 4
$ws.Range("E17").Value = 100
$ws.Range("F17").Value = 30
$ws.Range("H17").Value = -26.829268292682
$ws.Range("I17").Value = 67
$ws.Range("J17").Value = 65
$ws.Range("K17").Value = 3.076923076923
$ws.Range("L17").Value = -4.285714285714
$ws.Range("M17").Value = 91.428571428571
$ws.Range("N17").Value = 21.818181818181
# Row 18
$ws.Range("D18").Value = 6
$ws.Range("E18").Value = -66.666666666666
$ws.Range("F18").Value = 17
$ws.Range("G18").Value = 13
$ws.Range("H18").Value = 30.769230769230
$ws.Range("I18").Value = 23
$ws.Range("J18").Value = 35
$ws.Range("K18").Value = -34.285714285714
$ws.Range("L18").Value = 43.75
$ws.Range("M18").Value = -48.888888888888
$ws.Range("N18").Value = -94.147582697201
# Row 19
$ws.Range("C19").Value = 11
$ws.Range("D19").Value = 17
$ws.Range("E19").Value = -35.294117647058
$ws.Range("F19").Value = 59
$ws.Range("G19").Value = 73
$ws.Range("H19").Value = -19.178082191780
$ws.Range("I19").Value = 100
$ws.Range("J19").Value = 157
$ws.Range("K19").Value = -36.305732484076
$ws.Range("L19").Value = -8.256880733944
$ws.Range("M19").Value = 51.515151515151
$ws.Range("N19").Value = -44.134078212290
# Row 20
$ws.Range("C20").Value = 5
$ws.Range("E20").Value = 150
$ws.Range("G20").Value = 14
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 25
$ws.Range("J20").Value = 35
$ws.Range("K20").Value = -28.571428571428
$ws.Range("L20").Value = -47.916666666666
$ws.Range("M20").Value = -26.470588235294
$ws.Range("N20").Value = -92.236024844720
# Row 21
$ws.Range("C21").Value = 35
$ws.Range("D21").Value = 37
$ws.Range("E21").Value = -5.405405405405
$ws.Range("F21").Value = 142
$ws.Range("G21").Value = 169
$ws.Range("H21").Value = -15.976331360946
$ws.Range("I21").Value = 256
$ws.Range("J21").Value = 348
$ws.Range("K21").Value = -26.436781609195
$ws.Range("L21").Value = -9.540636042402
$ws.Range("M21").Value = 9.871244635193
$ws.Range("N21").Value = -77.969018932874
# Row 22
$ws.Range("C22").Value = 1
$ws.Range("G22").Value = 4
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 9
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = -10
$ws.Range("M22").Value = 80
# Row 24
$ws.Range("C24").Value = 26
$ws.Range("D24").Value = 63
$ws.Range("E24").Value = -58.730158730158
$ws.Range("F24").Value = 113
$ws.Range("G24").Value = 213
$ws.Range("H24").Value = -46.948356807511
$ws.Range("I24").Value = 205
$ws.Range("J24").Value = 379
$ws.Range("K24").Value = -45.910290237467
$ws.Range("L24").Value = -36.923076923076
$ws.Range("M24").Value = 49.635036496350
# Row 25
$ws.Range("C25").Value = 11
$ws.Range("D25").Value = 36
$ws.Range("E25").Value = -69.444444444444
$ws.Range("F25").Value = 57
$ws.Range("G25").Value = 128
$ws.Range("H25").Value = -55.46875
$ws.Range("I25").Value = 96
$ws.Range("J25").Value = 230
$ws.Range("K25").Value = -58.260869565217
$ws.Range("L25").Value = -49.738219895288
# Row 26
$ws.Range("C26").Value = 9
$ws.Range("D26").Value = 17
$ws.Range("E26").Value = -47.058823529411
$ws.Range("F26").Value = 68
$ws.Range("G26").Value = 86
$ws.Range("H26").Value = -20.930232558139
$ws.Range("I26").Value = 147
$ws.Range("J26").Value = 148
$ws.Range("K26").Value = -0.675675675675
$ws.Range("L26").Value = 17.6
$ws.Range("M26").Value = 14.84375
# Row 27
$ws.Range("C14").Copy($ws.Range("C27"))
$ws.Range("E27").Value = -100
$ws.Range("F27").Value = 3
$ws.Range("G27").Value = 3
$ws.Range("H27").Value = 0
$ws.Range("J27").Value = 10
$ws.Range("K27").Value = -20
$ws.Range("L27").Value = 0
# Row 28
$ws.Range("D28").Value = 3
$ws.Range("E28").Value = -33.333333333333
$ws.Range("F28").Value = 6
$ws.Range("G28").Value = 15
$ws.Range("H28").Value = -60
$ws.Range("I28").Value = 12
$ws.Range("J28").Value = 21
$ws.Range("K28").Value = -42.857142857142
$ws.Range("L28").Value = -7.692307692307
# Row 31
$ws.Range("G31").Value = 2
